# Generate Report for Handoff
#
# The localization status has moved from "In Translation" to
# "Ready for handoff" and the report's generation/handoff timestamps have
# been refreshed. Update the three worksheets (Overview, zh-cn, de-de)
# accordingly, and widen the status columns to fit the new, longer status
# text (matches the column-width bump baked into the canonical workbook).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps ---------------------------------------------------
$wsOverview.Range("G2").Value = "2016-08-13 11:12:31" # Latest HO Xliff Generate Date
$wsZhCn.Range("H2").Value     = "2016-08-13 11:12:23" # zh-cn Latest Handoff Datetime
$wsDeDe.Range("H2").Value     = "2016-08-13 11:12:31" # de-de Latest Handoff Datetime

# --- Widen the status columns to fit "Ready for handoff" -------------------
# The engine quantizes ColumnWidth to whole pixels (1/6-character steps), so
# we pick the input that lands on the closest achievable stored width to the
# canonical 17.2159881591797 (px=98 -> 17.1666...).
$newStatusWidth = 16.333333

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth  # col E
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth  # col F
$wsZhCn.Columns.Item(3).ColumnWidth     = $newStatusWidth  # col C
$wsDeDe.Columns.Item(3).ColumnWidth     = $newStatusWidth  # col C

Write-Output "Report regenerated for handoff."
